$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New odds data for rows 2-6, columns F (6) through AO (41)
$data = @(
    @(1.22,1.23,26,30,6.6,7.2,0,0,5,1.23,1.91,2.06,1.25,4.8,2.28,1.64,1.03,4.8,1000,1000,1000,1000,5.2,9.4,1000,1000,4.8,8.6,990,1000,9.6,980,1000,1000,16,1000),
    @(4.7,4.9,2.8,2.92,2.22,2.28,0,1.04,1.24,1.31,1.29,3.6,1.01,60,2.32,1.42,1.53,1.26,2.98,4.5,19.5,110,7,11,42,280,44,65,240,1000,350,1000,1000,1000,1000,610),
    @(1.16,410,1.05,980,3.35,48,1.46,1.04,3,1.41,1.7,2.32,1.25,4.7,1.03,1.01,1.15,2.08,980,1000,1000,1000,1000,990,990,1000,1000,990,990,1000,1000,1000,1000,1000,1000,1000),
    @(2.42,2.62,3.25,3.65,3.15,3.5,1.6,1.11,2.88,1.47,1.62,2.3,1.22,5,2.06,1.86,1.34,1.62,11.5,13.5,30,95,8.4,7.6,15.5,55,18,12.5,22,75,46,42,60,200,36,80),
    @(1.7,1.76,6,6.6,3.7,4.1,1.45,1.09,3.55,1.38,1.82,2.1,1.32,3.9,2,1.86,1.19,2.32,15,19,55,1000,7.8,9,25,480,9.800000000000001,10.5,42,1000,17.5,19.5,80,1000,12.5,1000)
)

$startRow = 2
$startCol = 6  # Column F

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $c = $startCol + $j
        $ws.Cells.Item($r, $c).Value = $rowValues[$j]
    }
}
